# "processing fticrms data from TMP1"
#
# The sample list had an extra, duplicated header row (row 45 repeated
# the "Vial File Position" / "Sample_ID" header that already lives in
# row 1) sitting in the middle of the data. Remove that stray row so
# the list is contiguous again (this also shifts every row below it up
# by one, e.g. old row 46 -> new row 45, ... old row 65 -> new row 64).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(45).Delete()

# Leave the view where the author ended up while reviewing the cleaned
# up list: scrolled down with D50 as the active cell.
$excel.ActiveWindow.ScrollRow = 19
$ws.Range("D50").Select() | Out-Null
